$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line contingencies (line7, line8) are inserted right after line6,
# pushing the former extr1..extr8 rows down by two rows. extr7 and extr8
# land on two brand-new rows (16 and 17) at the bottom of the table.
# Clone the formatting of the last existing data row onto the two new rows
# before filling in their values, so the new rows match the existing
# (bold + bordered) row style used for column A.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

$data = @(
    @(8,  "line7", 14, 11, $true),
    @(9,  "line8", 16, 9,  $false),
    @(10, "extr1", 5,  12, $true),
    @(11, "extr2", 5,  9,  $true),
    @(12, "extr3", 10, 11, $true),
    @(13, "extr4", 7,  8,  $true),
    @(14, "extr5", 9,  11, $true),
    @(15, "extr6", 7,  11, $true),
    @(16, "extr7", 5,  7,  $true),
    @(17, "extr8", 8,  5,  $true)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $r - 2
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
